$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.0467633
$ws.Range("C2").Value = 0.0487335
$ws.Range("B3").Value = 0.07045319999999999
$ws.Range("C3").Value = 0.06695479999999999
$ws.Range("B4").Value = 0.06944699999999999
$ws.Range("C4").Value = 0.073643
$ws.Range("B5").Value = 0.0687446
$ws.Range("C5").Value = 0.0684003
$ws.Range("B6").Value = 0.07184649999999999
$ws.Range("C6").Value = 0.0763292
$ws.Range("B7").Value = 0.0677532
$ws.Range("C7").Value = 0.0765948
$ws.Range("B8").Value = 0.06634710000000001
$ws.Range("C8").Value = 0.07066119999999999
$ws.Range("B9").Value = 0.06604549999999999
$ws.Range("C9").Value = 0.0701191
$ws.Range("B10").Value = 0.0651026
$ws.Range("C10").Value = 0.0713932
$ws.Range("B11").Value = 0.06485009999999999
$ws.Range("C11").Value = 0.0757284
$ws.Range("B12").Value = 0.06514739999999999
$ws.Range("C12").Value = 0.0735261
$ws.Range("B13").Value = 0.0643561
$ws.Range("C13").Value = 0.0734255
$ws.Range("B14").Value = 0.06486889999999999
$ws.Range("C14").Value = 0.0746039
$ws.Range("B15").Value = 0.06350889999999999
$ws.Range("C15").Value = 0.07324840000000001
$ws.Range("B16").Value = 0.06327149999999999
$ws.Range("C16").Value = 0.07385799999999999
$ws.Range("B17").Value = 0.0625245
$ws.Range("C17").Value = 0.07481309999999999
$ws.Range("B18").Value = 0.0614723
$ws.Range("C18").Value = 0.0752569
$ws.Range("B19").Value = 0.0613532
$ws.Range("C19").Value = 0.0757472
$ws.Range("B20").Value = 0.0616428
$ws.Range("C20").Value = 0.07676769999999999
$ws.Range("B21").Value = 0.0639245
$ws.Range("C21").Value = 0.077359
$ws.Range("B22").Value = 0.0624081
$ws.Range("C22").Value = 0.07936559999999999
$ws.Range("B23").Value = 0.0597885
$ws.Range("C23").Value = 0.079052
$ws.Range("B24").Value = 0.0584506
$ws.Range("C24").Value = 0.0815068
$ws.Range("B25").Value = 0.062557
$ws.Range("C25").Value = 0.0803644
$ws.Range("B26").Value = 0.0642065
$ws.Range("C26").Value = 0.08063380000000001
$ws.Range("B27").Value = 0.0566824
$ws.Range("C27").Value = 0.0796277
$ws.Range("B28").Value = 0.0569589
$ws.Range("C28").Value = 0.0804494
$ws.Range("B29").Value = 0.0556992
$ws.Range("C29").Value = 0.08298659999999999
$ws.Range("B30").Value = 0.0549723
$ws.Range("C30").Value = 0.08191229999999999
$ws.Range("B31").Value = 0.0544828
$ws.Range("C31").Value = 0.0835911
$ws.Range("B32").Value = 0.0549093
$ws.Range("C32").Value = 0.0839536
$ws.Range("B33").Value = 0.0535964
$ws.Range("C33").Value = 0.0850148
$ws.Range("B34").Value = 0.0537806
$ws.Range("C34").Value = 0.0836634
$ws.Range("B35").Value = 0.0540613
$ws.Range("C35").Value = 0.08407050000000001
$ws.Range("B36").Value = 0.0535874
$ws.Range("C36").Value = 0.08498219999999999
$ws.Range("B37").Value = 0.0522843
$ws.Range("C37").Value = 0.0859414
$ws.Range("B38").Value = 0.0519608
$ws.Range("C38").Value = 0.0864109
$ws.Range("B39").Value = 0.050965
$ws.Range("C39").Value = 0.0864138
$ws.Range("B40").Value = 0.0508956
$ws.Range("C40").Value = 0.08933720000000001
$ws.Range("B41").Value = 0.0534747
$ws.Range("C41").Value = 0.08889619999999999
$ws.Range("B42").Value = 0.0684694
$ws.Range("C42").Value = 0.0975973
$ws.Range("B43").Value = 0.0504215
$ws.Range("C43").Value = 0.0925718
$ws.Range("B44").Value = 0.0493211
$ws.Range("C44").Value = 0.09393899999999999
$ws.Range("B45").Value = 0.0481868
$ws.Range("C45").Value = 0.09102789999999999
$ws.Range("B46").Value = 0.0476639
$ws.Range("C46").Value = 0.0913003
$ws.Range("B47").Value = 0.0497213
$ws.Range("C47").Value = 0.0909706
$ws.Range("B48").Value = 0.0460525
$ws.Range("C48").Value = 0.0980834
$ws.Range("B49").Value = 0.0501845
$ws.Range("C49").Value = 0.0922331
$ws.Range("B50").Value = 0.0455816
$ws.Range("C50").Value = 0.09217500000000001
$ws.Range("B51").Value = 0.0459533
$ws.Range("C51").Value = 0.0967678
$ws.Range("B52").Value = 0.0442078
$ws.Range("C52").Value = 0.0943138
$ws.Range("B53").Value = 0.0459689
$ws.Range("C53").Value = 0.09400650000000001
$ws.Range("B54").Value = 0.0435796
$ws.Range("C54").Value = 0.0962325
$ws.Range("B55").Value = 0.043721
$ws.Range("C55").Value = 0.0948382
$ws.Range("B56").Value = 0.0428141
$ws.Range("C56").Value = 0.09578979999999999
$ws.Range("B57").Value = 0.0420995
$ws.Range("C57").Value = 0.100265
$ws.Range("B58").Value = 0.0419807
$ws.Range("C58").Value = 0.0973691
$ws.Range("B59").Value = 0.0417504
$ws.Range("C59").Value = 0.0975861
$ws.Range("B60").Value = 0.0416881
$ws.Range("C60").Value = 0.09812800000000001
$ws.Range("B61").Value = 0.0403392
$ws.Range("C61").Value = 0.0991388
$ws.Range("B62").Value = 0.039616
$ws.Range("C62").Value = 0.0999492
$ws.Range("B63").Value = 0.0389807
$ws.Range("C63").Value = 0.103598
$ws.Range("B64").Value = 0.0386416
$ws.Range("C64").Value = 0.100801
$ws.Range("B65").Value = 0.0407718
$ws.Range("C65").Value = 0.101056
$ws.Range("B66").Value = 0.0374552
$ws.Range("C66").Value = 0.101664
$ws.Range("B67").Value = 0.0367964
$ws.Range("C67").Value = 0.101588
$ws.Range("B68").Value = 0.0366471
$ws.Range("C68").Value = 0.102827
$ws.Range("B69").Value = 0.035756
$ws.Range("C69").Value = 0.105723
$ws.Range("B70").Value = 0.0358421
$ws.Range("C70").Value = 0.104518
$ws.Range("B71").Value = 0.0361845
$ws.Range("C71").Value = 0.103434
$ws.Range("B72").Value = 0.0341179
$ws.Range("C72").Value = 0.105703
$ws.Range("B73").Value = 0.0335061
$ws.Range("C73").Value = 0.104974
$ws.Range("B74").Value = 0.0352373
$ws.Range("C74").Value = 0.111499
$ws.Range("B75").Value = 0.0329565
$ws.Range("C75").Value = 0.11336
$ws.Range("B76").Value = 0.0332007
$ws.Range("C76").Value = 0.107513
$ws.Range("B77").Value = 0.0324623
$ws.Range("C77").Value = 0.112839
$ws.Range("B78").Value = 0.0314621
$ws.Range("C78").Value = 0.10879
$ws.Range("B79").Value = 0.0315148
$ws.Range("C79").Value = 0.111756
$ws.Range("B80").Value = 0.0302888
$ws.Range("C80").Value = 0.144872
$ws.Range("B81").Value = 0.032548
$ws.Range("C81").Value = 0.114621
$ws.Range("B82").Value = 0.0298351
$ws.Range("C82").Value = 0.111347
$ws.Range("B83").Value = 0.029173
$ws.Range("C83").Value = 0.110588
$ws.Range("B84").Value = 0.0285799
$ws.Range("C84").Value = 0.111092
$ws.Range("B85").Value = 0.0282174
$ws.Range("C85").Value = 0.115505
$ws.Range("B86").Value = 0.0280146
$ws.Range("C86").Value = 0.111858
$ws.Range("B87").Value = 0.0271155
$ws.Range("C87").Value = 0.113191
$ws.Range("B88").Value = 0.026867
$ws.Range("C88").Value = 0.113801
$ws.Range("B89").Value = 0.0265065
$ws.Range("C89").Value = 0.114367
$ws.Range("B90").Value = 0.0254739
$ws.Range("C90").Value = 0.114017
$ws.Range("B91").Value = 0.0250202
$ws.Range("C91").Value = 0.114645
$ws.Range("B92").Value = 0.0243298
$ws.Range("C92").Value = 0.119862
$ws.Range("B93").Value = 0.0243701
$ws.Range("C93").Value = 0.116532
$ws.Range("B94").Value = 0.0242332
$ws.Range("C94").Value = 0.116866
$ws.Range("B95").Value = 0.0231271
$ws.Range("C95").Value = 0.125036
$ws.Range("B96").Value = 0.0224824
$ws.Range("C96").Value = 0.117722
$ws.Range("B97").Value = 0.0219052
$ws.Range("C97").Value = 0.122868
$ws.Range("B98").Value = 0.0219857
$ws.Range("C98").Value = 0.118383
$ws.Range("B99").Value = 0.0213325
$ws.Range("C99").Value = 0.120636
$ws.Range("B100").Value = 0.0204447
$ws.Range("C100").Value = 0.127449
$ws.Range("B101").Value = 0.0204218
$ws.Range("C101").Value = 0.120909
$ws.Range("B102").Value = 0.0203905
$ws.Range("C102").Value = 0.120248
